{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// prompt in the document body with its updated value from the commit.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"2025-09-01 Monday\", \"2025-09-02 Tuesday\"],\n  [\"81\u00d760=\", \"67\u00d732=\"],\n  [\"29\u00d736=\", \"66\u00d746=\"],\n  [\"91\u00d773=\", \"66\u00d746=\"],\n  [\"37\u00d747=\", \"25\u00d765=\"],\n  [\"53\u00d746=\", \"22\u00d723=\"],\n  [\"79\u00d790=\", \"19\u00d780=\"],\n  [\"29\u00d794=\", \"49\u00d781=\"],\n  [\"92\u00d786=\", \"81\u00d782=\"],\n  [\"84\u00d783=\", \"39\u00d722=\"],\n  [\"76\u00d761=\", \"21\u00d787=\"],\n  [\"81\u00d777=\", \"40\u00d794=\"],\n  [\"44\u00d790=\", \"58\u00d788=\"],\n  [\"40\u00d735=\", \"99\u00d782=\"],\n  [\"34\u00d795=\", \"36\u00d772=\"],\n  [\"21\u00d721=\", \"63\u00d739=\"],\n  [\"50\u00d716=\", \"39\u00d717=\"],\n  [\"20\u00d723=\", \"46\u00d793=\"],\n  [\"70\u00d781=\", \"27\u00d794=\"],\n  [\"31\u00d712=\", \"30\u00d748=\"],\n  [\"83\u00d722=\", \"48\u00d746=\"],\n  [\"69\u00d796=\", \"16\u00d718=\"],\n  [\"31\u00d788=\", \"37\u00d721=\"],\n  [\"30\u00d785=\", \"84\u00d775=\"],\n  [\"77\u00d738=\", \"65\u00d772=\"],\n  [\"72\u00d792=\", \"75\u00d798=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit-by-two-digit multiplication\n# prompt in the document body with its updated value from the commit.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $result) {\n    Write-Output (\"WARNING: replace failed for \" + $oldText)\n  }\n}\n\nReplace-Text \"2025-09-01 Monday\" \"2025-09-02 Tuesday\"\nReplace-Text \"81\u00d760=\" \"67\u00d732=\"\nReplace-Text \"29\u00d736=\" \"66\u00d746=\"\nReplace-Text \"91\u00d773=\" \"66\u00d746=\"\nReplace-Text \"37\u00d747=\" \"25\u00d765=\"\nReplace-Text \"53\u00d746=\" \"22\u00d723=\"\nReplace-Text \"79\u00d790=\" \"19\u00d780=\"\nReplace-Text \"29\u00d794=\" \"49\u00d781=\"\nReplace-Text \"92\u00d786=\" \"81\u00d782=\"\nReplace-Text \"84\u00d783=\" \"39\u00d722=\"\nReplace-Text \"76\u00d761=\" \"21\u00d787=\"\nReplace-Text \"81\u00d777=\" \"40\u00d794=\"\nReplace-Text \"44\u00d790=\" \"58\u00d788=\"\nReplace-Text \"40\u00d735=\" \"99\u00d782=\"\nReplace-Text \"34\u00d795=\" \"36\u00d772=\"\nReplace-Text \"21\u00d721=\" \"63\u00d739=\"\nReplace-Text \"50\u00d716=\" \"39\u00d717=\"\nReplace-Text \"20\u00d723=\" \"46\u00d793=\"\nReplace-Text \"70\u00d781=\" \"27\u00d794=\"\nReplace-Text \"31\u00d712=\" \"30\u00d748=\"\nReplace-Text \"83\u00d722=\" \"48\u00d746=\"\nReplace-Text \"69\u00d796=\" \"16\u00d718=\"\nReplace-Text \"31\u00d788=\" \"37\u00d721=\"\nReplace-Text \"30\u00d785=\" \"84\u00d775=\"\nReplace-Text \"77\u00d738=\" \"65\u00d772=\"\nReplace-Text \"72\u00d792=\" \"75\u00d798=\"\n"}
